# Included objects for ActionRequired and ActionOverdue Count.
#
# On the "Objects_Navigation" sheet, insert two new rows just above the
# "User Site Sub Menu Transmittals Navigation" section (old row 10) to hold
# two new object-repository entries used to read the Action Required /
# Actions Overdue counters from the usersite menu. Inserting the rows shifts
# everything below down by two, which Excel keeps consistent for merged
# cells and data validations automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Navigation")

# Insert two blank rows at row 9 (pushes old row 9+ down to row 11+).
$null = $ws.Rows("9:10").Insert()

# New row 9: Action Required count locator.
$ws.Range("E9").Value = ".//*[@class='action-number action-required-number']"
$ws.Range("B9").Value = "Usersite Menu - Action Required Count"
$ws.Range("C9").Value = "xpath"
$ws.Range("D9").Value = "textbox_gettext"

# New row 10: Actions Overdue count locator.
$ws.Range("B10").Value = "Usersite Menu - Actions Overdue Count"
$ws.Range("C10").Value = "xpath"
$ws.Range("D10").Value = "textbox_gettext"
$ws.Range("E10").Value = ".//*[@class='action-number action-overdue-number']"

# Match the author's final selection position on the sheet.
$null = $ws.Range("D5").Select()
